# Adds a new "2022-Q4" quarter sheet (cloned from the "2022-Q3" sheet so it
# inherits identical number formats/styles), fills it with the Q4 fund
# holdings data, and inserts a corresponding summary row at the top of the
# "总计" (totals) sheet, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q4" worksheet by duplicating "2022-Q3" (keeps
#    fonts/borders/number-formats identical to its sibling quarter sheets)
#    and placing it immediately before "2022-Q3".
# ---------------------------------------------------------------------
$q3Sheet = $wb.Worksheets.Item("2022-Q3")
$q3Sheet.Copy($q3Sheet)
$newSheet = $wb.Worksheets.Item("2022-Q3 (2)")
$newSheet.Name = "2022-Q4"

# The template sheet has 12 data rows (2..13); the new sheet only needs 3
# (2..4), so drop the extra rows and shift the rest up.
$newSheet.Range("A5:H13").Delete(-4162)

# ---------------------------------------------------------------------
# 2. Populate the Q4 fund table. Column A (index) already reads 0,1,2 from
#    the template and needs no change. Columns B/D/E/F/G hold numeric-
#    looking values that must stay TEXT (matching every sibling sheet), so
#    they are written with a leading apostrophe to force text entry.
#    Column C (fund name) is ordinary text; column H (rank) is numeric.
# ---------------------------------------------------------------------
$newSheet.Range("B2").Value = "'001707"
$newSheet.Range("C2").Value = "诺安高端制造股票A"
$newSheet.Range("D2").Value = "'1.14"
$newSheet.Range("E2").Value = "'92.36"
$newSheet.Range("F2").Value = "'2.71"
$newSheet.Range("G2").Value = "'0.0309"
$newSheet.Range("H2").Value = 7

$newSheet.Range("B3").Value = "'000965"
$newSheet.Range("C3").Value = "汇丰晋信新动力混合"
$newSheet.Range("D3").Value = "'0.95"
$newSheet.Range("E3").Value = "'91.04"
$newSheet.Range("F3").Value = "'2.82"
$newSheet.Range("G3").Value = "'0.0268"
$newSheet.Range("H3").Value = 10

$newSheet.Range("B4").Value = "'014536"
$newSheet.Range("C4").Value = "诺安高端制造股票C"
$newSheet.Range("D4").Value = "'0.01"
$newSheet.Range("E4").Value = "'92.36"
$newSheet.Range("F4").Value = "'2.71"
$newSheet.Range("G4").Value = "'0.0003"
$newSheet.Range("H4").Value = 7

# ---------------------------------------------------------------------
# 3. Insert the matching summary row into the "总计" sheet: push the
#    existing 2022-Q3..2020-Q4 rows down by one (bottom-up copy so sources
#    aren't clobbered before being read) and write the new 2022-Q4 row on
#    top.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

for ($r = 9; $r -ge 2; $r--) {
    $dateVal = $totalSheet.Cells.Item($r, 2).Value()
    $countVal = $totalSheet.Cells.Item($r, 3).Value()
    $mvVal = $totalSheet.Cells.Item($r, 4).Value()

    $totalSheet.Cells.Item($r + 1, 2).Value = $dateVal
    $totalSheet.Cells.Item($r + 1, 3).Value = $countVal
    $totalSheet.Cells.Item($r + 1, 4).Value = $mvVal
}

# Column A is just the zero-based row index and already holds 0..7 in rows
# 2..9; only the brand-new row 10 needs a value (copy A9's style first so
# the new cell matches the rest of the column).
$totalSheet.Range("A9").Copy()
$totalSheet.Range("A10").PasteSpecial(-4122)
$totalSheet.Range("A10").Value = 8

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 0.06

# ---------------------------------------------------------------------
# 4. Restore the originally-active tab ("2020-Q4") since adding/copying
#    sheets moves the selection to the newly created sheet.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
